$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value to set across columns J:AS
$updates = @{
    100 = 351066.4746
    101 = 789501.6767
    102 = 187696.9708
    103 = 931753.1528
    104 = 270051.1343
    105 = 472839.5758
    106 = 20646.66679
    107 = 1610316.465
    114 = 10941.26985
    115 = 15351643.43
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $rangeAddress = "J" + $row + ":AS" + $row
    $ws.Range($rangeAddress).Value = $value
}
